$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new time-log entry for 11/03 (date serial 45233):
# Task: Completed daily operations, 8 hours
$ws.Range("A19").Value = 45233
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B19").Value = "Internship"

$ws.Range("C19").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Move the active selection down to the next empty row, as Excel does
# after the user finishes typing the row.
$ws.Range("C20").Select()
